# Motor_disability.docx edit
# 1) "Das größte Problem ist die Präzision..." paragraph:
#    - drop the direct w:sz="24" character/paragraph formatting
#    - re-split the runs and move the _GoBack bookmark so it wraps
#      "Das größte Proble|m| ist die ... Maus" (ending just before the
#      final period, which becomes its own trailing run)
# 2) "Tasten(schläge) werden visuell dargestellt" bullet:
#    - split around "schläge" and wrap it in spell-check proofErr tags

$d = $word.ActiveDocument

# --- Edit 1: locate the paragraph via the existing _GoBack bookmark ---
$bm = $d.Bookmarks.Item("_GoBack")
$p1 = $bm.Range.Paragraphs.First

$p1Xml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:r><w:t>Das größte</w:t></w:r><w:r><w:t xml:space="preserve"> Proble</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t xml:space="preserve"> ist die Präzision und Koordination von Bewegungen, vor allem mit den Händen mit einer Maus</w:t></w:r><w:bookmarkEnd w:id="0"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p1.Range.InsertXML($p1Xml) | Out-Null

# --- Edit 2: locate the "Tasten(schläge) ..." paragraph by its text ---
$p2 = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.StartsWith("Tasten(schl")) {
        $p2 = $para
    }
}

$p2Xml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Tasten(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>schläge</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>) werden visuell dargestellt</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$p2.Range.InsertXML($p2Xml) | Out-Null

Write-Output "edits applied"
